# Daily update at 8 AM UTC
# The sheet tracks a running win count per day. Each day the previous
# "latest" row (formatted as a plain date, no time) reverts to the
# standard date+time format, and a new row is appended for the new day
# carrying the "latest" date-only format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59 (2025-05-22) is no longer the most recent day -> give it back
# the regular date/time number format used by every earlier row.
$ws.Range("A59").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 60 (2025-05-23).
$ws.Range("A60").Value = 45800
$ws.Range("A60").NumberFormat = "YYYY-MM-DD"
$ws.Range("B60").Value = 250
$ws.Range("C60").Value = 261
$ws.Range("D60").Value = 251
